# Import dataroom structure: add a "Number"/"File" header row and shift the
# existing B:C table one column to the left (now A:B).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data currently lives in columns B (numbering) and C (names), with
# column A empty. Deleting the empty column A shifts B->A and C->B, freeing
# up row 1 above the data (which is still in rows 2-10) for a header.
$ws.Columns("A").Delete() | Out-Null

# Add the header row. Set B1 ("File") before A1 ("Number") so the new shared
# strings are appended in that order.
$ws.Range("B1").Value = "File"
$ws.Range("A1").Value = "Number"

# Select the whole populated table, matching the saved selection state.
$ws.Range("A1:B10").Select() | Out-Null
